$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.48559999999999
$ws.Range("B4").Value = 8.706799999999998
$ws.Range("B6").Value = 6.7406
$ws.Range("B7").Value = 5.048799999999999
$ws.Range("C7").Value = -14.05209999999999
$ws.Range("B8").Value = 6.630600000000006
$ws.Range("C11").Value = -12.19380000000001
$ws.Range("C12").Value = -11.2158
$ws.Range("D12").Value = -7.377399999999997
$ws.Range("E12").Value = 17.07640000000001
$ws.Range("D13").Value = -8.876700000000007
$ws.Range("E13").Value = 16.587
$ws.Range("D14").Value = -7.915099999999998
$ws.Range("C15").Value = -14.78029999999999
$ws.Range("B16").Value = 7.264799999999996
$ws.Range("D16").Value = -8.711300000000001
$ws.Range("D19").Value = -7.779599999999995
$ws.Range("B20").Value = 9.440399999999993
$ws.Range("C20").Value = -12.30510000000001
$ws.Range("D20").Value = -7.832999999999998
$ws.Range("B21").Value = 9.8277
$ws.Range("C21").Value = -12.33890000000001
$ws.Range("C22").Value = -12.2754
$ws.Range("D22").Value = -7.640400000000003
$ws.Range("E22").Value = 16.6112
$ws.Range("C23").Value = -12.0236
$ws.Range("E25").Value = 17.03380000000001
$ws.Range("B28").Value = 5.919399999999996
$ws.Range("B29").Value = 5.192100000000005
$ws.Range("C29").Value = -11.0389
$ws.Range("E29").Value = 17.43270000000001
$ws.Range("B30").Value = 4.954000000000001
$ws.Range("B32").Value = 7.338999999999998
$ws.Range("C34").Value = -12.51600000000001
$ws.Range("E34").Value = 17.5077
$ws.Range("D36").Value = -8.428000000000001
$ws.Range("B40").Value = 9.779599999999993
$ws.Range("C42").Value = -12.2786
$ws.Range("C43").Value = -13.57529999999999
$ws.Range("D43").Value = -8.086699999999997
$ws.Range("E43").Value = 16.3514
$ws.Range("C44").Value = -13.9128
$ws.Range("C45").Value = -13.6933
$ws.Range("B46").Value = 6.159399999999998
$ws.Range("C46").Value = -13.0778
$ws.Range("D46").Value = -8.502699999999997
$ws.Range("E48").Value = 17.43180000000001
$ws.Range("C50").Value = -13.92949999999999
$ws.Range("D50").Value = -8.279099999999998
$ws.Range("B51").Value = 6.412000000000002
$ws.Range("C51").Value = -12.06480000000001
$ws.Range("B52").Value = 5.409699999999999
$ws.Range("B57").Value = 6.306899999999996
$ws.Range("C57").Value = -13.29449999999999
$ws.Range("B59").Value = 5.427099999999998
$ws.Range("E60").Value = 15.9884
$ws.Range("B62").Value = 5.6293
$ws.Range("C65").Value = -12.62629999999999
$ws.Range("B66").Value = 5.683400000000002
$ws.Range("C66").Value = -11.62300000000001
$ws.Range("C67").Value = -11.6095
$ws.Range("E68").Value = 17.89050000000001
$ws.Range("E70").Value = 18.25360000000002
$ws.Range("E71").Value = 17.3517
$ws.Range("B73").Value = 8.468500000000001
$ws.Range("E73").Value = 17.66540000000001
$ws.Range("B74").Value = 9.336199999999989
$ws.Range("D76").Value = -7.986000000000003
$ws.Range("B77").Value = 8.729400000000004
$ws.Range("E78").Value = 17.06430000000002
$ws.Range("C79").Value = -11.9931
$ws.Range("C84").Value = -13.04279999999999
$ws.Range("C87").Value = -13.68829999999999
$ws.Range("E87").Value = 16.12869999999999
$ws.Range("B92").Value = 5.819299999999994
$ws.Range("C92").Value = -11.6186
$ws.Range("E92").Value = 18.21790000000002
$ws.Range("D95").Value = -7.993800000000002
$ws.Range("C97").Value = -12.3575
$ws.Range("D97").Value = -8.988400000000004
$ws.Range("D99").Value = -8.179200000000002
$ws.Range("B100").Value = 5.425099999999995
$ws.Range("E101").Value = 17.04300000000001
